$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 205; $r++) {
    $ws.Cells.Item($r, 3).Value = "2023-09-19"
}
